# Update the "Test Suite" worksheet so that every test case in column C
# (Runmode) is set to run ("Y"), reflecting "Running all the test cases".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Set the Runmode column (C2:C7) to "Y" for every row.
$ws.Range("C2:C7").Value = "Y"

# Update the active selection to match the edited range.
$ws.Range("C2:C7").Select()
